# Batch data excel added
# ---------------------------------------------------------------------------
# 1. Update the existing "loginData" sheet: C4 changes from "Feb@2023" (with a
#    style) to plain "Feb@2025".
# 2. Add a new "Batch" worksheet (after "loginData") with sample batch data.
# 3. Leave "Batch" as the active sheet / tab and restore the selection on the
#    "loginData" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. loginData sheet tweak
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C4").Value = "Feb@2025"
$ws1.Range("C4").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Add the "Batch" worksheet right after "loginData"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Batch"

# Header row
$ws2.Range("A1").Value = "scenario"
$ws2.Range("B1").Value = "ProgramName"
$ws2.Range("C1").Value = "BatchPrefix"
$ws2.Range("D1").Value = "BatchSuffix"
$ws2.Range("E1").Value = "Description"
$ws2.Range("F1").Value = "Status"
$ws2.Range("G1").Value = "NumberOfClasses"

# Row 2
$ws2.Range("A2").Value = "invalidSuffix"
$ws2.Range("D2").Value = "defg"

# Row 3
$ws2.Range("A3").Value = "invalidPrefix"
$ws2.Range("C3").Value = "tref"

# Row 4
$ws2.Range("A4").Value = "mandatoryFields"
$ws2.Range("B4").Value = "TestDef"
$ws2.Range("C4").Value = 4567
$ws2.Range("D4").Value = 6384
$ws2.Range("E4").Value = "testCase"
$ws2.Range("F4").Value = "ACTIVE"
$ws2.Range("G4").Value = 2

# Row 5
$ws2.Range("A5").Value = "oneFieldBlank"
$ws2.Range("B5").Value = "TestDef"
$ws2.Range("C5").Value = 4567
$ws2.Range("E5").Value = "testCase"
$ws2.Range("F5").Value = "ACTIVE"
$ws2.Range("G5").Value = 2

# Row 6 (styled cell: Consolas 10pt black)
$ws2.Range("A6").Value = "Invalid"
$ws2.Range("A6").Font.Name = "Consolas"
$ws2.Range("A6").Font.Size = 10
$ws2.Range("A6").Font.Color = 0

$ws2.Range("E6").Value = "r"
$ws2.Range("G6").Value = 0

# Row 7
$ws2.Range("A7").Value = "valid"
$ws2.Range("E7").Value = "testcases"
$ws2.Range("G7").Value = 1

# Column widths for the new sheet
$ws2.Columns.Item(2).ColumnWidth = 15.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 12.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 11.0
$ws2.Columns.Item(5).ColumnWidth = 14.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 16.666666666666664

# ---------------------------------------------------------------------------
# 3. Selections / active sheet
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A2:A8").Select()

$ws2.Activate()
$ws2.Range("D9").Select()
